# Make example numbers more realistic.
#
# The original workbook used inflated placeholder magnitudes (e.g. dividend
# amounts in the millions, production/cost flows in the tens of millions)
# and a couple of override fractions that didn't reflect the intended
# scenario. This trims those numbers down to more plausible values, flips
# one override flag on, and leaves the final selection on the "metric
# overrides" sheet (matching where the author ended up working).

$wb = $excel.ActiveWorkbook

# --- Sheet "metric overrides" ---
$wsOverrides = $wb.Worksheets.Item("metric overrides")
[void]($wsOverrides.Range("C4").Value = $true)

$wsOverrides.Range("C8").Value = 0.03
$wsOverrides.Range("E8").Value = 0.2
$wsOverrides.Range("F8").Value = 0.03

$wsOverrides.Range("C9").Value = 0.02
$wsOverrides.Range("E9").Value = 0.2
$wsOverrides.Range("F9").Value = 0.03

$wsOverrides.Range("C10").Value = 0.01
$wsOverrides.Range("E10").Value = 0.2
$wsOverrides.Range("F10").Value = 0.03

# --- Sheet "dividend" ---
$wsDividend = $wb.Worksheets.Item("dividend")
$wsDividend.Range("E5").Value = -1000
$wsDividend.Range("E6").Value = 1000

# --- Sheet "production" ---
$wsProduction = $wb.Worksheets.Item("production")
$wsProduction.Range("C3").Value = 10000
$wsProduction.Range("C4").Value = -100000
$wsProduction.Range("C5").Value = -200000

# --- Sheet "costs" ---
$wsCosts = $wb.Worksheets.Item("costs")
$wsCosts.Range("C3").Value = -10
$wsCosts.Range("C4").Value = -100
$wsCosts.Range("C5").Value = 10000

# --- Update selections / active sheet to match the saved view ---
[void]$wsDividend.Activate()
[void]$wsDividend.Range("E7").Select()

[void]$wsProduction.Activate()
[void]$wsProduction.Range("C9").Select()

[void]$wsCosts.Activate()
[void]$wsCosts.Range("C6").Select()

[void]$wsOverrides.Activate()
[void]$wsOverrides.Range("E11").Select()
